$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 41,22
$arr[0,0] = 75
$arr[0,1] = 'memory'
$arr[0,2] = 6
$arr[0,3] = 1
$arr[0,4] = 1
$arr[0,5] = 284
$arr[0,6] = 'living_rooms'
$arr[0,7] = 'living_rooms'
$arr[0,8] = $null
$arr[0,9] = 'new'
$arr[0,10] = 'f'
$arr[0,11] = 'stimuli/img_a9he3.png'
$arr[0,12] = 83.06521739130434
$arr[0,13] = 63.95652173913044
$arr[0,14] = 73.51086956521739
$arr[0,15] = 46
$arr[0,16] = 8
$arr[0,17] = 8
$arr[0,18] = 8
$arr[0,19] = 8
$arr[0,20] = 8
$arr[0,21] = 8
$arr[1,0] = 75
$arr[1,1] = 'memory'
$arr[1,2] = 6
$arr[1,3] = 1
$arr[1,4] = 2
$arr[1,5] = 285
$arr[1,6] = 'living_rooms'
$arr[1,7] = 'living_rooms'
$arr[1,8] = $null
$arr[1,9] = 'new'
$arr[1,10] = 'f'
$arr[1,11] = 'stimuli/img_o30wb.png'
$arr[1,12] = 81.06666666666666
$arr[1,13] = 65.37777777777778
$arr[1,14] = 73.22222222222223
$arr[1,15] = 45
$arr[1,16] = 8
$arr[1,17] = 8
$arr[1,18] = 8
$arr[1,19] = 8
$arr[1,20] = 8
$arr[1,21] = 8
$arr[2,0] = 75
$arr[2,1] = 'memory'
$arr[2,2] = 6
$arr[2,3] = 1
$arr[2,4] = 3
$arr[2,5] = 286
$arr[2,6] = 'living_rooms'
$arr[2,7] = 'living_rooms'
$arr[2,8] = $null
$arr[2,9] = 'new'
$arr[2,10] = 'f'
$arr[2,11] = 'stimuli/img_zxvl3.png'
$arr[2,12] = 68.78260869565217
$arr[2,13] = 47.56521739130435
$arr[2,14] = 58.17391304347827
$arr[2,15] = 46
$arr[2,16] = 5
$arr[2,17] = 5
$arr[2,18] = 5
$arr[2,19] = 5
$arr[2,20] = 5
$arr[2,21] = 5
$arr[3,0] = 75
$arr[3,1] = 'memory'
$arr[3,2] = 6
$arr[3,3] = 1
$arr[3,4] = 4
$arr[3,5] = 287
$arr[3,6] = 'living_rooms'
$arr[3,7] = 'living_rooms'
$arr[3,8] = 'target'
$arr[3,9] = 'old'
$arr[3,10] = 'j'
$arr[3,11] = 'stimuli/img_wbws6.png'
$arr[3,12] = 57.97777777777777
$arr[3,13] = 42.53333333333333
$arr[3,14] = 50.25555555555555
$arr[3,15] = 45
$arr[3,16] = 4
$arr[3,17] = 4
$arr[3,18] = 4
$arr[3,19] = 4
$arr[3,20] = 4
$arr[3,21] = 5
$arr[4,0] = 75
$arr[4,1] = 'memory'
$arr[4,2] = 6
$arr[4,3] = 1
$arr[4,4] = 5
$arr[4,5] = 288
$arr[4,6] = 'living_rooms'
$arr[4,7] = 'living_rooms'
$arr[4,8] = $null
$arr[4,9] = 'new'
$arr[4,10] = 'f'
$arr[4,11] = 'stimuli/img_5jy9c.png'
$arr[4,12] = 87.37209302325581
$arr[4,13] = 79.18604651162791
$arr[4,14] = 83.27906976744185
$arr[4,15] = 43
$arr[4,16] = 10
$arr[4,17] = 10
$arr[4,18] = 10
$arr[4,19] = 10
$arr[4,20] = 9
$arr[4,21] = 10
$arr[5,0] = 75
$arr[5,1] = 'memory'
$arr[5,2] = 6
$arr[5,3] = 1
$arr[5,4] = 6
$arr[5,5] = 289
$arr[5,6] = 'living_rooms'
$arr[5,7] = 'living_rooms'
$arr[5,8] = 'target'
$arr[5,9] = 'old'
$arr[5,10] = 'j'
$arr[5,11] = 'stimuli/img_pdzf1.png'
$arr[5,12] = 86.23913043478261
$arr[5,13] = 67.17391304347827
$arr[5,14] = 76.70652173913044
$arr[5,15] = 46
$arr[5,16] = 9
$arr[5,17] = 9
$arr[5,18] = 9
$arr[5,19] = 9
$arr[5,20] = 9
$arr[5,21] = 8
$arr[6,0] = 75
$arr[6,1] = 'memory'
$arr[6,2] = 6
$arr[6,3] = 1
$arr[6,4] = 7
$arr[6,5] = 290
$arr[6,6] = 'living_rooms'
$arr[6,7] = 'living_rooms'
$arr[6,8] = $null
$arr[6,9] = 'new'
$arr[6,10] = 'f'
$arr[6,11] = 'stimuli/img_9bkl9.png'
$arr[6,12] = 46.62162162162162
$arr[6,13] = 34.27027027027027
$arr[6,14] = 40.44594594594595
$arr[6,15] = 37
$arr[6,16] = 3
$arr[6,17] = 3
$arr[6,18] = 3
$arr[6,19] = 3
$arr[6,20] = 3
$arr[6,21] = 3
$arr[7,0] = 75
$arr[7,1] = 'memory'
$arr[7,2] = 6
$arr[7,3] = 1
$arr[7,4] = 8
$arr[7,5] = 291
$arr[7,6] = 'living_rooms'
$arr[7,7] = 'living_rooms'
$arr[7,8] = $null
$arr[7,9] = 'new'
$arr[7,10] = 'f'
$arr[7,11] = 'stimuli/img_g13d5.png'
$arr[7,12] = 73
$arr[7,13] = 51.51111111111111
$arr[7,14] = 62.25555555555556
$arr[7,15] = 45
$arr[7,16] = 6
$arr[7,17] = 6
$arr[7,18] = 6
$arr[7,19] = 6
$arr[7,20] = 6
$arr[7,21] = 6
$arr[8,0] = 75
$arr[8,1] = 'memory'
$arr[8,2] = 6
$arr[8,3] = 1
$arr[8,4] = 9
$arr[8,5] = 292
$arr[8,6] = 'living_rooms'
$arr[8,7] = 'living_rooms'
$arr[8,8] = 'target'
$arr[8,9] = 'old'
$arr[8,10] = 'j'
$arr[8,11] = 'stimuli/img_hmmra.png'
$arr[8,12] = 54.65853658536585
$arr[8,13] = 34.24390243902439
$arr[8,14] = 44.45121951219512
$arr[8,15] = 41
$arr[8,16] = 3
$arr[8,17] = 3
$arr[8,18] = 3
$arr[8,19] = 4
$arr[8,20] = 4
$arr[8,21] = 3
$arr[9,0] = 75
$arr[9,1] = 'memory'
$arr[9,2] = 6
$arr[9,3] = 1
$arr[9,4] = 10
$arr[9,5] = 293
$arr[9,6] = 'living_rooms'
$arr[9,7] = 'living_rooms'
$arr[9,8] = $null
$arr[9,9] = 'new'
$arr[9,10] = 'f'
$arr[9,11] = 'stimuli/img_165pk.png'
$arr[9,12] = 85.73333333333333
$arr[9,13] = 69.22222222222223
$arr[9,14] = 77.47777777777779
$arr[9,15] = 45
$arr[9,16] = 9
$arr[9,17] = 9
$arr[9,18] = 9
$arr[9,19] = 9
$arr[9,20] = 9
$arr[9,21] = 9
$arr[10,0] = 75
$arr[10,1] = 'memory'
$arr[10,2] = 6
$arr[10,3] = 1
$arr[10,4] = 11
$arr[10,5] = 294
$arr[10,6] = 'living_rooms'
$arr[10,7] = 'living_rooms'
$arr[10,8] = 'target'
$arr[10,9] = 'old'
$arr[10,10] = 'j'
$arr[10,11] = 'stimuli/img_5jp4f.png'
$arr[10,12] = 84.85714285714286
$arr[10,13] = 67.83333333333333
$arr[10,14] = 76.3452380952381
$arr[10,15] = 42
$arr[10,16] = 9
$arr[10,17] = 9
$arr[10,18] = 9
$arr[10,19] = 8
$arr[10,20] = 8
$arr[10,21] = 9
$arr[11,0] = 75
$arr[11,1] = 'memory'
$arr[11,2] = 6
$arr[11,3] = 1
$arr[11,4] = 12
$arr[11,5] = 295
$arr[11,6] = 'living_rooms'
$arr[11,7] = 'living_rooms'
$arr[11,8] = 'target'
$arr[11,9] = 'old'
$arr[11,10] = 'j'
$arr[11,11] = 'stimuli/img_b21d7.png'
$arr[11,12] = 27.75555555555556
$arr[11,13] = 13.86666666666667
$arr[11,14] = 20.81111111111111
$arr[11,15] = 45
$arr[11,16] = 1
$arr[11,17] = 1
$arr[11,18] = 1
$arr[11,19] = 2
$arr[11,20] = 2
$arr[11,21] = 1
$arr[12,0] = 75
$arr[12,1] = 'memory'
$arr[12,2] = 6
$arr[12,3] = 1
$arr[12,4] = 13
$arr[12,5] = 296
$arr[12,6] = 'living_rooms'
$arr[12,7] = 'living_rooms'
$arr[12,8] = $null
$arr[12,9] = 'new'
$arr[12,10] = 'f'
$arr[12,11] = 'stimuli/img_jpldg.png'
$arr[12,12] = 79.54545454545455
$arr[12,13] = 57.75
$arr[12,14] = 68.64772727272728
$arr[12,15] = 44
$arr[12,16] = 7
$arr[12,17] = 7
$arr[12,18] = 7
$arr[12,19] = 7
$arr[12,20] = 7
$arr[12,21] = 7
$arr[13,0] = 75
$arr[13,1] = 'memory'
$arr[13,2] = 6
$arr[13,3] = 1
$arr[13,4] = 14
$arr[13,5] = 297
$arr[13,6] = 'living_rooms'
$arr[13,7] = 'living_rooms'
$arr[13,8] = 'target'
$arr[13,9] = 'old'
$arr[13,10] = 'j'
$arr[13,11] = 'stimuli/img_tn8ys.png'
$arr[13,12] = 86.70454545454545
$arr[13,13] = 72.4090909090909
$arr[13,14] = 79.55681818181819
$arr[13,15] = 44
$arr[13,16] = 10
$arr[13,17] = 10
$arr[13,18] = 10
$arr[13,19] = 9
$arr[13,20] = 9
$arr[13,21] = 10
$arr[14,0] = 75
$arr[14,1] = 'memory'
$arr[14,2] = 6
$arr[14,3] = 1
$arr[14,4] = 15
$arr[14,5] = 298
$arr[14,6] = 'living_rooms'
$arr[14,7] = 'living_rooms'
$arr[14,8] = $null
$arr[14,9] = 'new'
$arr[14,10] = 'f'
$arr[14,11] = 'stimuli/img_3sw8t.png'
$arr[14,12] = 67.4888888888889
$arr[14,13] = 48.51111111111111
$arr[14,14] = 58
$arr[14,15] = 45
$arr[14,16] = 5
$arr[14,17] = 5
$arr[14,18] = 5
$arr[14,19] = 5
$arr[14,20] = 5
$arr[14,21] = 5
$arr[15,0] = 75
$arr[15,1] = 'memory'
$arr[15,2] = 6
$arr[15,3] = 1
$arr[15,4] = 16
$arr[15,5] = 299
$arr[15,6] = 'living_rooms'
$arr[15,7] = 'living_rooms'
$arr[15,8] = 'target'
$arr[15,9] = 'old'
$arr[15,10] = 'j'
$arr[15,11] = 'stimuli/img_16kib.png'
$arr[15,12] = 80.97727272727273
$arr[15,13] = 61.11363636363637
$arr[15,14] = 71.04545454545455
$arr[15,15] = 44
$arr[15,16] = 8
$arr[15,17] = 8
$arr[15,18] = 8
$arr[15,19] = 7
$arr[15,20] = 7
$arr[15,21] = 7
$arr[16,0] = 75
$arr[16,1] = 'memory'
$arr[16,2] = 6
$arr[16,3] = 1
$arr[16,4] = 17
$arr[16,5] = 300
$arr[16,6] = 'living_rooms'
$arr[16,7] = 'living_rooms'
$arr[16,8] = $null
$arr[16,9] = 'new'
$arr[16,10] = 'f'
$arr[16,11] = 'stimuli/img_pbsj1.png'
$arr[16,12] = 73.88636363636364
$arr[16,13] = 51.52272727272727
$arr[16,14] = 62.70454545454545
$arr[16,15] = 44
$arr[16,16] = 6
$arr[16,17] = 6
$arr[16,18] = 6
$arr[16,19] = 6
$arr[16,20] = 6
$arr[16,21] = 6
$arr[17,0] = 75
$arr[17,1] = 'memory'
$arr[17,2] = 6
$arr[17,3] = 1
$arr[17,4] = 18
$arr[17,5] = 301
$arr[17,6] = 'living_rooms'
$arr[17,7] = 'living_rooms'
$arr[17,8] = 'target'
$arr[17,9] = 'old'
$arr[17,10] = 'j'
$arr[17,11] = 'stimuli/img_gka64.png'
$arr[17,12] = 19.23809523809524
$arr[17,13] = 20.02380952380953
$arr[17,14] = 19.63095238095238
$arr[17,15] = 42
$arr[17,16] = 1
$arr[17,17] = 1
$arr[17,18] = 1
$arr[17,19] = 1
$arr[17,20] = 1
$arr[17,21] = 2
$arr[18,0] = 75
$arr[18,1] = 'memory'
$arr[18,2] = 6
$arr[18,3] = 1
$arr[18,4] = 19
$arr[18,5] = 302
$arr[18,6] = 'living_rooms'
$arr[18,7] = 'living_rooms'
$arr[18,8] = $null
$arr[18,9] = 'new'
$arr[18,10] = 'f'
$arr[18,11] = 'stimuli/img_i6wsx.png'
$arr[18,12] = 79.07142857142857
$arr[18,13] = 58
$arr[18,14] = 68.53571428571428
$arr[18,15] = 42
$arr[18,16] = 7
$arr[18,17] = 7
$arr[18,18] = 7
$arr[18,19] = 7
$arr[18,20] = 7
$arr[18,21] = 7
$arr[19,0] = 75
$arr[19,1] = 'memory'
$arr[19,2] = 6
$arr[19,3] = 1
$arr[19,4] = 20
$arr[19,5] = 303
$arr[19,6] = 'living_rooms'
$arr[19,7] = 'living_rooms'
$arr[19,8] = 'target'
$arr[19,9] = 'old'
$arr[19,10] = 'j'
$arr[19,11] = 'stimuli/img_xr3up.png'
$arr[19,12] = 76.24444444444444
$arr[19,13] = 55.88888888888889
$arr[19,14] = 66.06666666666666
$arr[19,15] = 45
$arr[19,16] = 7
$arr[19,17] = 7
$arr[19,18] = 7
$arr[19,19] = 6
$arr[19,20] = 6
$arr[19,21] = 6
$arr[20,0] = 75
$arr[20,1] = 'memory'
$arr[20,2] = 6
$arr[20,3] = 1
$arr[20,4] = 21
$arr[20,5] = 304
$arr[20,6] = 'living_rooms'
$arr[20,7] = 'living_rooms'
$arr[20,8] = $null
$arr[20,9] = 'new'
$arr[20,10] = 'f'
$arr[20,11] = 'stimuli/img_xzyzy.png'
$arr[20,12] = 85.37209302325581
$arr[20,13] = 68.90697674418605
$arr[20,14] = 77.13953488372093
$arr[20,15] = 43
$arr[20,16] = 9
$arr[20,17] = 9
$arr[20,18] = 9
$arr[20,19] = 9
$arr[20,20] = 9
$arr[20,21] = 9
$arr[21,0] = 75
$arr[21,1] = 'memory'
$arr[21,2] = 6
$arr[21,3] = 1
$arr[21,4] = 22
$arr[21,5] = 305
$arr[21,6] = 'living_rooms'
$arr[21,7] = 'living_rooms'
$arr[21,8] = 'target'
$arr[21,9] = 'old'
$arr[21,10] = 'j'
$arr[21,11] = 'stimuli/img_c89x3.png'
$arr[21,12] = 72.8695652173913
$arr[21,13] = 49.65217391304348
$arr[21,14] = 61.26086956521739
$arr[21,15] = 46
$arr[21,16] = 6
$arr[21,17] = 6
$arr[21,18] = 6
$arr[21,19] = 6
$arr[21,20] = 6
$arr[21,21] = 5
$arr[22,0] = 75
$arr[22,1] = 'memory'
$arr[22,2] = 6
$arr[22,3] = 1
$arr[22,4] = 23
$arr[22,5] = 306
$arr[22,6] = 'living_rooms'
$arr[22,7] = 'living_rooms'
$arr[22,8] = 'target'
$arr[22,9] = 'old'
$arr[22,10] = 'j'
$arr[22,11] = 'stimuli/img_pjfx6.png'
$arr[22,12] = 32.23404255319149
$arr[22,13] = 26.59574468085106
$arr[22,14] = 29.41489361702127
$arr[22,15] = 47
$arr[22,16] = 2
$arr[22,17] = 2
$arr[22,18] = 2
$arr[22,19] = 2
$arr[22,20] = 2
$arr[22,21] = 3
$arr[23,0] = 75
$arr[23,1] = 'memory'
$arr[23,2] = 6
$arr[23,3] = 1
$arr[23,4] = 24
$arr[23,5] = 307
$arr[23,6] = 'living_rooms'
$arr[23,7] = 'living_rooms'
$arr[23,8] = 'target'
$arr[23,9] = 'old'
$arr[23,10] = 'j'
$arr[23,11] = 'stimuli/img_vgh2g.png'
$arr[23,12] = 93.81395348837209
$arr[23,13] = 78.27906976744185
$arr[23,14] = 86.04651162790697
$arr[23,15] = 43
$arr[23,16] = 10
$arr[23,17] = 10
$arr[23,18] = 10
$arr[23,19] = 10
$arr[23,20] = 10
$arr[23,21] = 10
$arr[24,0] = 75
$arr[24,1] = 'memory'
$arr[24,2] = 6
$arr[24,3] = 1
$arr[24,4] = 25
$arr[24,5] = 308
$arr[24,6] = 'living_rooms'
$arr[24,7] = 'living_rooms'
$arr[24,8] = $null
$arr[24,9] = 'new'
$arr[24,10] = 'f'
$arr[24,11] = 'stimuli/img_c0vzo.png'
$arr[24,12] = 21.51162790697675
$arr[24,13] = 8.232558139534884
$arr[24,14] = 14.87209302325581
$arr[24,15] = 43
$arr[24,16] = 1
$arr[24,17] = 1
$arr[24,18] = 1
$arr[24,19] = 1
$arr[24,20] = 1
$arr[24,21] = 1
$arr[25,0] = 75
$arr[25,1] = 'memory'
$arr[25,2] = 6
$arr[25,3] = 1
$arr[25,4] = 26
$arr[25,5] = 309
$arr[25,6] = 'living_rooms'
$arr[25,7] = 'living_rooms'
$arr[25,8] = $null
$arr[25,9] = 'new'
$arr[25,10] = 'f'
$arr[25,11] = 'stimuli/img_wgddx.png'
$arr[25,12] = 45.6304347826087
$arr[25,13] = 34.30434782608695
$arr[25,14] = 39.96739130434783
$arr[25,15] = 46
$arr[25,16] = 3
$arr[25,17] = 3
$arr[25,18] = 3
$arr[25,19] = 3
$arr[25,20] = 3
$arr[25,21] = 4
$arr[26,0] = 75
$arr[26,1] = 'memory'
$arr[26,2] = 6
$arr[26,3] = 1
$arr[26,4] = 27
$arr[26,5] = 310
$arr[26,6] = 'living_rooms'
$arr[26,7] = 'living_rooms'
$arr[26,8] = 'target'
$arr[26,9] = 'old'
$arr[26,10] = 'j'
$arr[26,11] = 'stimuli/img_j856a.png'
$arr[26,12] = 38.225
$arr[26,13] = 25.875
$arr[26,14] = 32.05
$arr[26,15] = 40
$arr[26,16] = 2
$arr[26,17] = 2
$arr[26,18] = 2
$arr[26,19] = 3
$arr[26,20] = 3
$arr[26,21] = 2
$arr[27,0] = 75
$arr[27,1] = 'memory'
$arr[27,2] = 6
$arr[27,3] = 1
$arr[27,4] = 28
$arr[27,5] = 311
$arr[27,6] = 'living_rooms'
$arr[27,7] = $null
$arr[27,8] = $null
$arr[27,9] = 'catch'
$arr[27,10] = 'f'
$arr[27,11] = 'stimuli/catch_05.jpg'
$arr[27,12] = $null
$arr[27,13] = $null
$arr[27,14] = $null
$arr[27,15] = $null
$arr[27,16] = $null
$arr[27,17] = $null
$arr[27,18] = $null
$arr[27,19] = $null
$arr[27,20] = $null
$arr[27,21] = $null
$arr[28,0] = 75
$arr[28,1] = 'memory'
$arr[28,2] = 6
$arr[28,3] = 1
$arr[28,4] = 29
$arr[28,5] = 312
$arr[28,6] = 'living_rooms'
$arr[28,7] = 'living_rooms'
$arr[28,8] = 'target'
$arr[28,9] = 'old'
$arr[28,10] = 'j'
$arr[28,11] = 'stimuli/img_q9lab.png'
$arr[28,12] = 53.97560975609756
$arr[28,13] = 32.90243902439025
$arr[28,14] = 43.4390243902439
$arr[28,15] = 41
$arr[28,16] = 3
$arr[28,17] = 3
$arr[28,18] = 3
$arr[28,19] = 3
$arr[28,20] = 4
$arr[28,21] = 3
$arr[29,0] = 75
$arr[29,1] = 'memory'
$arr[29,2] = 6
$arr[29,3] = 1
$arr[29,4] = 30
$arr[29,5] = 313
$arr[29,6] = 'living_rooms'
$arr[29,7] = 'living_rooms'
$arr[29,8] = $null
$arr[29,9] = 'new'
$arr[29,10] = 'f'
$arr[29,11] = 'stimuli/img_w8yhd.png'
$arr[29,12] = 55.74418604651163
$arr[29,13] = 38.90697674418605
$arr[29,14] = 47.32558139534883
$arr[29,15] = 43
$arr[29,16] = 4
$arr[29,17] = 4
$arr[29,18] = 4
$arr[29,19] = 4
$arr[29,20] = 4
$arr[29,21] = 4
$arr[30,0] = 75
$arr[30,1] = 'memory'
$arr[30,2] = 6
$arr[30,3] = 1
$arr[30,4] = 31
$arr[30,5] = 314
$arr[30,6] = 'living_rooms'
$arr[30,7] = 'living_rooms'
$arr[30,8] = $null
$arr[30,9] = 'new'
$arr[30,10] = 'f'
$arr[30,11] = 'stimuli/img_8dmpq.png'
$arr[30,12] = 30.65909090909091
$arr[30,13] = 24.11363636363636
$arr[30,14] = 27.38636363636364
$arr[30,15] = 44
$arr[30,16] = 2
$arr[30,17] = 2
$arr[30,18] = 2
$arr[30,19] = 2
$arr[30,20] = 2
$arr[30,21] = 2
$arr[31,0] = 75
$arr[31,1] = 'memory'
$arr[31,2] = 6
$arr[31,3] = 1
$arr[31,4] = 32
$arr[31,5] = 315
$arr[31,6] = 'living_rooms'
$arr[31,7] = 'living_rooms'
$arr[31,8] = 'target'
$arr[31,9] = 'old'
$arr[31,10] = 'j'
$arr[31,11] = 'stimuli/img_vgaye.png'
$arr[31,12] = 80.33333333333333
$arr[31,13] = 64.57777777777778
$arr[31,14] = 72.45555555555555
$arr[31,15] = 45
$arr[31,16] = 8
$arr[31,17] = 8
$arr[31,18] = 8
$arr[31,19] = 8
$arr[31,20] = 7
$arr[31,21] = 8
$arr[32,0] = 75
$arr[32,1] = 'memory'
$arr[32,2] = 6
$arr[32,3] = 1
$arr[32,4] = 33
$arr[32,5] = 316
$arr[32,6] = 'living_rooms'
$arr[32,7] = 'living_rooms'
$arr[32,8] = $null
$arr[32,9] = 'new'
$arr[32,10] = 'f'
$arr[32,11] = 'stimuli/img_jkm86.png'
$arr[32,12] = 58.32558139534883
$arr[32,13] = 38.65116279069768
$arr[32,14] = 48.48837209302326
$arr[32,15] = 43
$arr[32,16] = 4
$arr[32,17] = 4
$arr[32,18] = 4
$arr[32,19] = 4
$arr[32,20] = 4
$arr[32,21] = 4
$arr[33,0] = 75
$arr[33,1] = 'memory'
$arr[33,2] = 6
$arr[33,3] = 1
$arr[33,4] = 34
$arr[33,5] = 317
$arr[33,6] = 'living_rooms'
$arr[33,7] = 'living_rooms'
$arr[33,8] = 'target'
$arr[33,9] = 'old'
$arr[33,10] = 'j'
$arr[33,11] = 'stimuli/img_6a0hu.png'
$arr[33,12] = 61.275
$arr[33,13] = 42.025
$arr[33,14] = 51.65
$arr[33,15] = 40
$arr[33,16] = 4
$arr[33,17] = 4
$arr[33,18] = 4
$arr[33,19] = 5
$arr[33,20] = 4
$arr[33,21] = 5
$arr[34,0] = 75
$arr[34,1] = 'memory'
$arr[34,2] = 6
$arr[34,3] = 1
$arr[34,4] = 35
$arr[34,5] = 318
$arr[34,6] = 'living_rooms'
$arr[34,7] = 'living_rooms'
$arr[34,8] = $null
$arr[34,9] = 'new'
$arr[34,10] = 'f'
$arr[34,11] = 'stimuli/img_rych7.png'
$arr[34,12] = 30.4468085106383
$arr[34,13] = 23.4468085106383
$arr[34,14] = 26.9468085106383
$arr[34,15] = 47
$arr[34,16] = 2
$arr[34,17] = 2
$arr[34,18] = 2
$arr[34,19] = 2
$arr[34,20] = 2
$arr[34,21] = 2
$arr[35,0] = 75
$arr[35,1] = 'memory'
$arr[35,2] = 6
$arr[35,3] = 1
$arr[35,4] = 36
$arr[35,5] = 319
$arr[35,6] = 'living_rooms'
$arr[35,7] = 'living_rooms'
$arr[35,8] = 'target'
$arr[35,9] = 'old'
$arr[35,10] = 'j'
$arr[35,11] = 'stimuli/img_ra2nm.png'
$arr[35,12] = 70.75
$arr[35,13] = 50.375
$arr[35,14] = 60.5625
$arr[35,15] = 40
$arr[35,16] = 6
$arr[35,17] = 6
$arr[35,18] = 6
$arr[35,19] = 5
$arr[35,20] = 5
$arr[35,21] = 6
$arr[36,0] = 75
$arr[36,1] = 'memory'
$arr[36,2] = 6
$arr[36,3] = 1
$arr[36,4] = 37
$arr[36,5] = 320
$arr[36,6] = 'living_rooms'
$arr[36,7] = 'living_rooms'
$arr[36,8] = $null
$arr[36,9] = 'new'
$arr[36,10] = 'f'
$arr[36,11] = 'stimuli/img_nb8p4.png'
$arr[36,12] = 16.36170212765957
$arr[36,13] = 12.70212765957447
$arr[36,14] = 14.53191489361702
$arr[36,15] = 47
$arr[36,16] = 1
$arr[36,17] = 1
$arr[36,18] = 1
$arr[36,19] = 1
$arr[36,20] = 1
$arr[36,21] = 1
$arr[37,0] = 75
$arr[37,1] = 'memory'
$arr[37,2] = 6
$arr[37,3] = 1
$arr[37,4] = 38
$arr[37,5] = 321
$arr[37,6] = 'living_rooms'
$arr[37,7] = 'living_rooms'
$arr[37,8] = 'target'
$arr[37,9] = 'old'
$arr[37,10] = 'j'
$arr[37,11] = 'stimuli/img_x4bln.png'
$arr[37,12] = 76.34042553191489
$arr[37,13] = 59.51063829787234
$arr[37,14] = 67.92553191489361
$arr[37,15] = 47
$arr[37,16] = 7
$arr[37,17] = 7
$arr[37,18] = 7
$arr[37,19] = 7
$arr[37,20] = 7
$arr[37,21] = 7
$arr[38,0] = 75
$arr[38,1] = 'memory'
$arr[38,2] = 6
$arr[38,3] = 1
$arr[38,4] = 39
$arr[38,5] = 322
$arr[38,6] = 'living_rooms'
$arr[38,7] = 'living_rooms'
$arr[38,8] = 'target'
$arr[38,9] = 'old'
$arr[38,10] = 'j'
$arr[38,11] = 'stimuli/img_z4jxm.png'
$arr[38,12] = 88.30952380952381
$arr[38,13] = 72.64285714285714
$arr[38,14] = 80.47619047619048
$arr[38,15] = 42
$arr[38,16] = 10
$arr[38,17] = 10
$arr[38,18] = 10
$arr[38,19] = 10
$arr[38,20] = 10
$arr[38,21] = 10
$arr[39,0] = 75
$arr[39,1] = 'memory'
$arr[39,2] = 6
$arr[39,3] = 1
$arr[39,4] = 40
$arr[39,5] = 323
$arr[39,6] = 'living_rooms'
$arr[39,7] = 'living_rooms'
$arr[39,8] = $null
$arr[39,9] = 'new'
$arr[39,10] = 'f'
$arr[39,11] = 'stimuli/img_dg5h7.png'
$arr[39,12] = 88.72093023255815
$arr[39,13] = 76.06976744186046
$arr[39,14] = 82.3953488372093
$arr[39,15] = 43
$arr[39,16] = 10
$arr[39,17] = 10
$arr[39,18] = 10
$arr[39,19] = 10
$arr[39,20] = 10
$arr[39,21] = 10
$arr[40,0] = 75
$arr[40,1] = 'memory'
$arr[40,2] = 6
$arr[40,3] = 1
$arr[40,4] = 41
$arr[40,5] = 324
$arr[40,6] = 'living_rooms'
$arr[40,7] = 'living_rooms'
$arr[40,8] = 'target'
$arr[40,9] = 'old'
$arr[40,10] = 'j'
$arr[40,11] = 'stimuli/img_j4ttn.png'
$arr[40,12] = 12.61904761904762
$arr[40,13] = 11.42857142857143
$arr[40,14] = 12.02380952380952
$arr[40,15] = 42
$arr[40,16] = 1
$arr[40,17] = 1
$arr[40,18] = 1
$arr[40,19] = 1
$arr[40,20] = 1
$arr[40,21] = 1

$ws.Range("A2:V42").Value = $arr

Write-Output "done"
